$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description cell in row 2 with the new spicule type text
$ws.Range("A2").Value = "no spicules"

# The previously existing (empty) cells B2:F2 are no longer present in the
# updated sheet, so clear their contents entirely.
$ws.Range("B2:F2").ClearContents()
